$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new attendance/clock-in log entry as row 4.
# Date (A), ClockIn (C), ClockOut (D) and Log (E) are plain text values
# (matching the existing rows, which are not real Excel dates/times);
# EmployeeID (B) is numeric. A leading apostrophe forces Excel to treat
# the date-looking string as text instead of auto-converting it to a
# date serial; resetting the style back to "Normal" afterwards drops the
# quote-prefix formatting so the new cells stay unstyled like rows 2-3.

$ws.Range("A4").Value = "'10/05/2025"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = 2702258535

$ws.Range("C4").Value = "14:06:02"

$ws.Range("D4").Value = "'"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = "'"
$ws.Range("E4").Style = "Normal"
